$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "ETH"
$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Service-account-mosip-resident-client"
